# Applies the "car details in progress" edit to Samochody.xlsx
# - Fixes the class (column D) of row 107 from "C" to "D"
# - Fills in columns E:T (engine/spec details) for rows 102-126
# - Updates the sheet view (selection / scroll position)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kompletne dane")

# Row 107 also had its class corrected from "C" to "D"
$ws.Cells.Item(107, 4).Value = "D"

# Row 102
$ws.Cells.Item(102, 5).Value = 4.5999999999999996
$ws.Cells.Item(102, 6).Value = 8
$ws.Cells.Item(102, 7).Value = 274
$ws.Cells.Item(102, 8).Value = 2
$ws.Cells.Item(102, 9).Value = 2
$ws.Cells.Item(102, 10).Value = 8
$ws.Cells.Item(102, 11).Value = 190
$ws.Cells.Item(102, 12).Value = 393
$ws.Cells.Item(102, 13).Value = 1292
$ws.Cells.Item(102, 14).Value = 339
$ws.Cells.Item(102, 15).Value = 90
$ws.Cells.Item(102, 16).Value = 450
$ws.Cells.Item(102, 17).Value = 185
$ws.Cells.Item(102, 18).Value = 132
$ws.Cells.Item(102, 19).Value = 259
$ws.Cells.Item(102, 20).Value = 4

# Row 103
$ws.Cells.Item(103, 5).Value = 3
$ws.Cells.Item(103, 6).Value = 12
$ws.Cells.Item(103, 7).Value = 300
$ws.Cells.Item(103, 8).Value = 2
$ws.Cells.Item(103, 9).Value = 2
$ws.Cells.Item(103, 10).Value = 6.2
$ws.Cells.Item(103, 11).Value = 265
$ws.Cells.Item(103, 12).Value = 343
$ws.Cells.Item(103, 13).Value = 1000
$ws.Cells.Item(103, 14).Value = 246
$ws.Cells.Item(103, 15).Value = 133
$ws.Cells.Item(103, 16).Value = 441
$ws.Cells.Item(103, 17).Value = 168
$ws.Cells.Item(103, 18).Value = 124
$ws.Cells.Item(103, 19).Value = 240
$ws.Cells.Item(103, 20).Value = 5

# Row 104
$ws.Cells.Item(104, 5).Value = 3.4
$ws.Cells.Item(104, 6).Value = 8
$ws.Cells.Item(104, 7).Value = 300
$ws.Cells.Item(104, 8).Value = 2
$ws.Cells.Item(104, 9).Value = 2
$ws.Cells.Item(104, 10).Value = 5.6
$ws.Cells.Item(104, 11).Value = 275
$ws.Cells.Item(104, 12).Value = 324
$ws.Cells.Item(104, 13).Value = 1393
$ws.Cells.Item(104, 14).Value = 200
$ws.Cells.Item(104, 15).Value = 95
$ws.Cells.Item(104, 16).Value = 423
$ws.Cells.Item(104, 17).Value = 189
$ws.Cells.Item(104, 18).Value = 117
$ws.Cells.Item(104, 19).Value = 245
$ws.Cells.Item(104, 20).Value = 5

# Row 105
$ws.Cells.Item(105, 5).Value = 4
$ws.Cells.Item(105, 6).Value = 8
$ws.Cells.Item(105, 7).Value = 500
$ws.Cells.Item(105, 8).Value = 2
$ws.Cells.Item(105, 9).Value = 2
$ws.Cells.Item(105, 10).Value = 2.8
$ws.Cells.Item(105, 11).Value = 300
$ws.Cells.Item(105, 12).Value = 500
$ws.Cells.Item(105, 13).Value = 1120
$ws.Cells.Item(105, 14).Value = 0
$ws.Cells.Item(105, 15).Value = 121
$ws.Cells.Item(105, 16).Value = 500
$ws.Cells.Item(105, 17).Value = 196
$ws.Cells.Item(105, 18).Value = 114
$ws.Cells.Item(105, 19).Value = 274
$ws.Cells.Item(105, 20).Value = 6

# Row 106
$ws.Cells.Item(106, 5).Value = 2.5
$ws.Cells.Item(106, 6).Value = 5
$ws.Cells.Item(106, 7).Value = 225
$ws.Cells.Item(106, 8).Value = 5
$ws.Cells.Item(106, 9).Value = 3
$ws.Cells.Item(106, 10).Value = 6.8
$ws.Cells.Item(106, 11).Value = 241
$ws.Cells.Item(106, 12).Value = 320
$ws.Cells.Item(106, 13).Value = 1392
$ws.Cells.Item(106, 14).Value = 385
$ws.Cells.Item(106, 15).Value = 55
$ws.Cells.Item(106, 16).Value = 436
$ws.Cells.Item(106, 17).Value = 184
$ws.Cells.Item(106, 18).Value = 145
$ws.Cells.Item(106, 19).Value = 264
$ws.Cells.Item(106, 20).Value = 6

# Row 107
$ws.Cells.Item(107, 5).Value = 2.2999999999999998
$ws.Cells.Item(107, 6).Value = 6
$ws.Cells.Item(107, 7).Value = 125
$ws.Cells.Item(107, 8).Value = 4
$ws.Cells.Item(107, 9).Value = 2
$ws.Cells.Item(107, 10).Value = 10
$ws.Cells.Item(107, 11).Value = 186
$ws.Cells.Item(107, 12).Value = 176
$ws.Cells.Item(107, 13).Value = 931
$ws.Cells.Item(107, 14).Value = 260
$ws.Cells.Item(107, 15).Value = 58
$ws.Cells.Item(107, 16).Value = 426
$ws.Cells.Item(107, 17).Value = 165
$ws.Cells.Item(107, 18).Value = 128
$ws.Cells.Item(107, 19).Value = 256
$ws.Cells.Item(107, 20).Value = 4

# Row 108
$ws.Cells.Item(108, 5).Value = 3.5
$ws.Cells.Item(108, 6).Value = 6
$ws.Cells.Item(108, 7).Value = 542
$ws.Cells.Item(108, 8).Value = 2
$ws.Cells.Item(108, 9).Value = 2
$ws.Cells.Item(108, 10).Value = 3.8
$ws.Cells.Item(108, 11).Value = 342
$ws.Cells.Item(108, 12).Value = 640
$ws.Cells.Item(108, 13).Value = 1350
$ws.Cells.Item(108, 14).Value = 110
$ws.Cells.Item(108, 15).Value = 120
$ws.Cells.Item(108, 16).Value = 493
$ws.Cells.Item(108, 17).Value = 200
$ws.Cells.Item(108, 18).Value = 115
$ws.Cells.Item(108, 19).Value = 264
$ws.Cells.Item(108, 20).Value = 5

# Row 109
$ws.Cells.Item(109, 5).Value = 2.2999999999999998
$ws.Cells.Item(109, 6).Value = 6
$ws.Cells.Item(109, 7).Value = 150
$ws.Cells.Item(109, 8).Value = 2
$ws.Cells.Item(109, 9).Value = 2
$ws.Cells.Item(109, 10).Value = 9.6999999999999993
$ws.Cells.Item(109, 11).Value = 200
$ws.Cells.Item(109, 12).Value = 196
$ws.Cells.Item(109, 13).Value = 1295
$ws.Cells.Item(109, 14).Value = 340
$ws.Cells.Item(109, 15).Value = 65
$ws.Cells.Item(109, 16).Value = 429
$ws.Cells.Item(109, 17).Value = 176
$ws.Cells.Item(109, 18).Value = 129
$ws.Cells.Item(109, 19).Value = 240
$ws.Cells.Item(109, 20).Value = 5

# Row 110
$ws.Cells.Item(110, 5).Value = 1.6
$ws.Cells.Item(110, 6).Value = 4
$ws.Cells.Item(110, 7).Value = 160
$ws.Cells.Item(110, 8).Value = 5
$ws.Cells.Item(110, 9).Value = 3
$ws.Cells.Item(110, 10).Value = 7.3
$ws.Cells.Item(110, 11).Value = 215
$ws.Cells.Item(110, 12).Value = 150
$ws.Cells.Item(110, 13).Value = 1080
$ws.Cells.Item(110, 14).Value = 190
$ws.Cells.Item(110, 15).Value = 45
$ws.Cells.Item(110, 16).Value = 408
$ws.Cells.Item(110, 17).Value = 170
$ws.Cells.Item(110, 18).Value = 135
$ws.Cells.Item(110, 19).Value = 257
$ws.Cells.Item(110, 20).Value = 5

# Row 111
$ws.Cells.Item(111, 5).Value = 4.8
$ws.Cells.Item(111, 6).Value = 12
$ws.Cells.Item(111, 7).Value = 375
$ws.Cells.Item(111, 8).Value = 2
$ws.Cells.Item(111, 9).Value = 2
$ws.Cells.Item(111, 10).Value = 5.0999999999999996
$ws.Cells.Item(111, 11).Value = 298
$ws.Cells.Item(111, 12).Value = 410
$ws.Cells.Item(111, 13).Value = 1480
$ws.Cells.Item(111, 14).Value = 240
$ws.Cells.Item(111, 15).Value = 120
$ws.Cells.Item(111, 16).Value = 414
$ws.Cells.Item(111, 17).Value = 200
$ws.Cells.Item(111, 18).Value = 107
$ws.Cells.Item(111, 19).Value = 245
$ws.Cells.Item(111, 20).Value = 5

# Row 112
$ws.Cells.Item(112, 5).Value = 4
$ws.Cells.Item(112, 6).Value = 8
$ws.Cells.Item(112, 7).Value = 560
$ws.Cells.Item(112, 8).Value = 5
$ws.Cells.Item(112, 9).Value = 5
$ws.Cells.Item(112, 10).Value = 3.9
$ws.Cells.Item(112, 11).Value = 305
$ws.Cells.Item(112, 12).Value = 700
$ws.Cells.Item(112, 13).Value = 1950
$ws.Cells.Item(112, 14).Value = 565
$ws.Cells.Item(112, 15).Value = 65
$ws.Cells.Item(112, 16).Value = 498
$ws.Cells.Item(112, 17).Value = 194
$ws.Cells.Item(112, 18).Value = 146
$ws.Cells.Item(112, 19).Value = 292
$ws.Cells.Item(112, 20).Value = 8

# Row 113
$ws.Cells.Item(113, 5).Value = 4.7
$ws.Cells.Item(113, 6).Value = 12
$ws.Cells.Item(113, 7).Value = 520
$ws.Cells.Item(113, 8).Value = 2
$ws.Cells.Item(113, 9).Value = 2
$ws.Cells.Item(113, 10).Value = 3.9
$ws.Cells.Item(113, 11).Value = 325
$ws.Cells.Item(113, 12).Value = 471
$ws.Cells.Item(113, 13).Value = 1230
$ws.Cells.Item(113, 14).Value = 57
$ws.Cells.Item(113, 15).Value = 105
$ws.Cells.Item(113, 16).Value = 448
$ws.Cells.Item(113, 17).Value = 199
$ws.Cells.Item(113, 18).Value = 112
$ws.Cells.Item(113, 19).Value = 258
$ws.Cells.Item(113, 20).Value = 6

# Row 114
$ws.Cells.Item(114, 5).Value = 3.4
$ws.Cells.Item(114, 6).Value = 8
$ws.Cells.Item(114, 7).Value = 767
$ws.Cells.Item(114, 8).Value = 2
$ws.Cells.Item(114, 9).Value = 2
$ws.Cells.Item(114, 10).Value = 3
$ws.Cells.Item(114, 11).Value = 322
$ws.Cells.Item(114, 12).Value = 1280
$ws.Cells.Item(114, 13).Value = 1490
$ws.Cells.Item(114, 14).Value = 110
$ws.Cells.Item(114, 15).Value = 70
$ws.Cells.Item(114, 16).Value = 464
$ws.Cells.Item(114, 17).Value = 194
$ws.Cells.Item(114, 18).Value = 117
$ws.Cells.Item(114, 19).Value = 273
$ws.Cells.Item(114, 20).Value = 7

# Row 115
$ws.Cells.Item(115, 5).Value = 5.2
$ws.Cells.Item(115, 6).Value = 10
$ws.Cells.Item(115, 7).Value = 570
$ws.Cells.Item(115, 8).Value = 2
$ws.Cells.Item(115, 9).Value = 2
$ws.Cells.Item(115, 10).Value = 3.4
$ws.Cells.Item(115, 11).Value = 320
$ws.Cells.Item(115, 12).Value = 540
$ws.Cells.Item(115, 13).Value = 1430
$ws.Cells.Item(115, 14).Value = 110
$ws.Cells.Item(115, 15).Value = 70
$ws.Cells.Item(115, 16).Value = 439
$ws.Cells.Item(115, 17).Value = 190
$ws.Cells.Item(115, 18).Value = 117
$ws.Cells.Item(115, 19).Value = 256
$ws.Cells.Item(115, 20).Value = 6

# Row 116
$ws.Cells.Item(116, 5).Value = 2
$ws.Cells.Item(116, 6).Value = 4
$ws.Cells.Item(116, 7).Value = 240
$ws.Cells.Item(116, 8).Value = 2
$ws.Cells.Item(116, 9).Value = 2
$ws.Cells.Item(116, 10).Value = 6.2
$ws.Cells.Item(116, 11).Value = 240
$ws.Cells.Item(116, 12).Value = 208
$ws.Cells.Item(116, 13).Value = 1270
$ws.Cells.Item(116, 14).Value = 143
$ws.Cells.Item(116, 15).Value = 50
$ws.Cells.Item(116, 16).Value = 414
$ws.Cells.Item(116, 17).Value = 175
$ws.Cells.Item(116, 18).Value = 127
$ws.Cells.Item(116, 19).Value = 240
$ws.Cells.Item(116, 20).Value = 6

# Row 117
$ws.Cells.Item(117, 5).Value = 4
$ws.Cells.Item(117, 6).Value = 8
$ws.Cells.Item(117, 7).Value = 720
$ws.Cells.Item(117, 8).Value = 2
$ws.Cells.Item(117, 9).Value = 2
$ws.Cells.Item(117, 10).Value = 2.9
$ws.Cells.Item(117, 11).Value = 341
$ws.Cells.Item(117, 12).Value = 770
$ws.Cells.Item(117, 13).Value = 1419
$ws.Cells.Item(117, 14).Value = 360
$ws.Cells.Item(117, 15).Value = 72
$ws.Cells.Item(117, 16).Value = 354
$ws.Cells.Item(117, 17).Value = 193
$ws.Cells.Item(117, 18).Value = 120
$ws.Cells.Item(117, 19).Value = 267
$ws.Cells.Item(117, 20).Value = 7

# Row 118
$ws.Cells.Item(118, 5).Value = 4
$ws.Cells.Item(118, 6).Value = 8
$ws.Cells.Item(118, 7).Value = 522
$ws.Cells.Item(118, 8).Value = 2
$ws.Cells.Item(118, 9).Value = 2
$ws.Cells.Item(118, 10).Value = 3.8
$ws.Cells.Item(118, 11).Value = 310
$ws.Cells.Item(118, 12).Value = 670
$ws.Cells.Item(118, 13).Value = 1570
$ws.Cells.Item(118, 14).Value = 285
$ws.Cells.Item(118, 15).Value = 65
$ws.Cells.Item(118, 16).Value = 454
$ws.Cells.Item(118, 17).Value = 207
$ws.Cells.Item(118, 18).Value = 129
$ws.Cells.Item(118, 19).Value = 263
$ws.Cells.Item(118, 20).Value = 7

# Row 119
$ws.Cells.Item(119, 5).Value = 3.9
$ws.Cells.Item(119, 6).Value = 8
$ws.Cells.Item(119, 7).Value = 670
$ws.Cells.Item(119, 8).Value = 2
$ws.Cells.Item(119, 9).Value = 2
$ws.Cells.Item(119, 10).Value = 3
$ws.Cells.Item(119, 11).Value = 325
$ws.Cells.Item(119, 12).Value = 760
$ws.Cells.Item(119, 13).Value = 1420
$ws.Cells.Item(119, 14).Value = 230
$ws.Cells.Item(119, 15).Value = 78
$ws.Cells.Item(119, 16).Value = 457
$ws.Cells.Item(119, 17).Value = 195
$ws.Cells.Item(119, 18).Value = 121
$ws.Cells.Item(119, 19).Value = 265
$ws.Cells.Item(119, 20).Value = 7

# Row 120
$ws.Cells.Item(120, 5).Value = 3.8
$ws.Cells.Item(120, 6).Value = 6
$ws.Cells.Item(120, 7).Value = 570
$ws.Cells.Item(120, 8).Value = 4
$ws.Cells.Item(120, 9).Value = 2
$ws.Cells.Item(120, 10).Value = 2.8
$ws.Cells.Item(120, 11).Value = 315
$ws.Cells.Item(120, 12).Value = 637
$ws.Cells.Item(120, 13).Value = 1752
$ws.Cells.Item(120, 14).Value = 315
$ws.Cells.Item(120, 15).Value = 74
$ws.Cells.Item(120, 16).Value = 471
$ws.Cells.Item(120, 17).Value = 190
$ws.Cells.Item(120, 18).Value = 137
$ws.Cells.Item(120, 19).Value = 178
$ws.Cells.Item(120, 20).Value = 6

# Row 121
$ws.Cells.Item(121, 5).Value = 5.2
$ws.Cells.Item(121, 6).Value = 12
$ws.Cells.Item(121, 7).Value = 725
$ws.Cells.Item(121, 8).Value = 4
$ws.Cells.Item(121, 9).Value = 2
$ws.Cells.Item(121, 10).Value = 3.4
$ws.Cells.Item(121, 11).Value = 340
$ws.Cells.Item(121, 12).Value = 900
$ws.Cells.Item(121, 13).Value = 1693
$ws.Cells.Item(121, 14).Value = 186
$ws.Cells.Item(121, 15).Value = 78
$ws.Cells.Item(121, 16).Value = 471
$ws.Cells.Item(121, 17).Value = 194
$ws.Cells.Item(121, 18).Value = 128
$ws.Cells.Item(121, 19).Value = 281
$ws.Cells.Item(121, 20).Value = 8

# Row 122
$ws.Cells.Item(122, 5).Value = 3
$ws.Cells.Item(122, 6).Value = 6
$ws.Cells.Item(122, 7).Value = 310
$ws.Cells.Item(122, 8).Value = 5
$ws.Cells.Item(122, 9).Value = 5
$ws.Cells.Item(122, 10).Value = 5.6
$ws.Cells.Item(122, 11).Value = 275
$ws.Cells.Item(122, 12).Value = 440
$ws.Cells.Item(122, 13).Value = 1785
$ws.Cells.Item(122, 14).Value = 535
$ws.Cells.Item(122, 15).Value = 65
$ws.Cells.Item(122, 16).Value = 497
$ws.Cells.Item(122, 17).Value = 191
$ws.Cells.Item(122, 18).Value = 142
$ws.Cells.Item(122, 19).Value = 291
$ws.Cells.Item(122, 20).Value = 7

# Row 123
$ws.Cells.Item(123, 5).Value = 1.8
$ws.Cells.Item(123, 6).Value = 4
$ws.Cells.Item(123, 7).Value = 192
$ws.Cells.Item(123, 8).Value = 4
$ws.Cells.Item(123, 9).Value = 2
$ws.Cells.Item(123, 10).Value = 7.4
$ws.Cells.Item(123, 11).Value = 225
$ws.Cells.Item(123, 12).Value = 180
$ws.Cells.Item(123, 13).Value = 1215
$ws.Cells.Item(123, 14).Value = 365
$ws.Cells.Item(123, 15).Value = 55
$ws.Cells.Item(123, 16).Value = 434
$ws.Cells.Item(123, 17).Value = 174
$ws.Cells.Item(123, 18).Value = 132
$ws.Cells.Item(123, 19).Value = 260
$ws.Cells.Item(123, 20).Value = 6

# Row 124
$ws.Cells.Item(124, 5).Value = 2
$ws.Cells.Item(124, 6).Value = 4
$ws.Cells.Item(124, 7).Value = 210
$ws.Cells.Item(124, 8).Value = 4
$ws.Cells.Item(124, 9).Value = 3
$ws.Cells.Item(124, 10).Value = 6.9
$ws.Cells.Item(124, 11).Value = 240
$ws.Cells.Item(124, 12).Value = 280
$ws.Cells.Item(124, 13).Value = 1298
$ws.Cells.Item(124, 14).Value = 312
$ws.Cells.Item(124, 15).Value = 55
$ws.Cells.Item(124, 16).Value = 426
$ws.Cells.Item(124, 17).Value = 181
$ws.Cells.Item(124, 18).Value = 140
$ws.Cells.Item(124, 19).Value = 258
$ws.Cells.Item(124, 20).Value = 6

# Row 125
$ws.Cells.Item(125, 5).Value = 3
$ws.Cells.Item(125, 6).Value = 6
$ws.Cells.Item(125, 7).Value = 591
$ws.Cells.Item(125, 8).Value = 4
$ws.Cells.Item(125, 9).Value = 2
$ws.Cells.Item(125, 10).Value = 4.0999999999999996
$ws.Cells.Item(125, 11).Value = 259
$ws.Cells.Item(125, 12).Value = 500
$ws.Cells.Item(125, 13).Value = 1540
$ws.Cells.Item(125, 14).Value = 0
$ws.Cells.Item(125, 15).Value = 51
$ws.Cells.Item(125, 16).Value = 438
$ws.Cells.Item(125, 17).Value = 185
$ws.Cells.Item(125, 18).Value = 129
$ws.Cells.Item(125, 19).Value = 247
$ws.Cells.Item(125, 20).Value = 6

# Row 126
$ws.Cells.Item(126, 5).Value = 4.9000000000000004
$ws.Cells.Item(126, 6).Value = 12
$ws.Cells.Item(126, 7).Value = 426
$ws.Cells.Item(126, 8).Value = 2
$ws.Cells.Item(126, 9).Value = 2
$ws.Cells.Item(126, 10).Value = 4.8
$ws.Cells.Item(126, 11).Value = 314
$ws.Cells.Item(126, 12).Value = 491
$ws.Cells.Item(126, 13).Value = 1590
$ws.Cells.Item(126, 14).Value = 249
$ws.Cells.Item(126, 15).Value = 100
$ws.Cells.Item(126, 16).Value = 448
$ws.Cells.Item(126, 17).Value = 198
$ws.Cells.Item(126, 18).Value = 114
$ws.Cells.Item(126, 19).Value = 255
$ws.Cells.Item(126, 20).Value = 5

# Restore the view state that was active when the edit was made
$ws.Range("J130").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 97
$win.ScrollColumn = 1

Write-Host "applied car detail rows 102-126"